$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.515.25"
$ws.Range("E2").Value = "  +3.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.985.95"
$ws.Range("E3").Value = "  +4.10%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.86"
$ws.Range("E5").Value = "  +6.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.88"
$ws.Range("E6").Value = "  +8.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.431"
$ws.Range("E8").Value = "  +7.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.47"
$ws.Range("E9").Value = "  +12.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  +12.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.353"
$ws.Range("E11").Value = "  +7.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.128"
$ws.Range("E12").Value = "  +4.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.496.42"
$ws.Range("E13").Value = "  +4.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.69"
$ws.Range("E14").Value = "  +12.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000153"
$ws.Range("E15").Value = "  +14.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.529.61"
$ws.Range("E16").Value = "  +4.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.986.47"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.80"
$ws.Range("E18").Value = "  +10.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.44"
$ws.Range("E19").Value = "  +8.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.79"
$ws.Range("E20").Value = "  +10.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.85"
$ws.Range("E21").Value = "  +7.43%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.474"
$ws.Range("E23").Value = "  +7.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.17"
$ws.Range("E24").Value = "  +5.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.994"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  +7.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0903"
$ws.Range("E27").Value = "  +12.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.44"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.94"
$ws.Range("E29").Value = "  +13.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.20"
$ws.Range("E30").Value = "  +7.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.76"
$ws.Range("E31").Value = "  +9.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.54"
$ws.Range("E32").Value = "  +8.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.64"
$ws.Range("E33").Value = "  +17.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("E34").Value = "  +6.88%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.56"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.26"
$ws.Range("E36").Value = "  +5.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0674"
$ws.Range("E37").Value = "  +10.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.24"
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.021.43"
$ws.Range("E39").Value = "  +4.37%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.34"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.643"
$ws.Range("E42").Value = "  +7.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.257.03"
$ws.Range("E43").Value = "  +10.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.40"
$ws.Range("E44").Value = "  +7.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.986"
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.58"
$ws.Range("E46").Value = "  +5.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.94"
$ws.Range("E47").Value = "  +22.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0236"
$ws.Range("E48").Value = "  +12.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.77"
$ws.Range("E49").Value = "  +9.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.11"
$ws.Range("E50").Value = "  +8.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0869"
$ws.Range("E51").Value = "  +11.00%  "
